$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above row 48; this shifts existing rows 48-110
# down to 49-111, matching the diff (old row N becomes new row N+1,
# and a brand new record is inserted at row 48).
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new weekly record.
$ws.Range("A48").Value = 2
$ws.Range("B48").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 44966
$ws.Range("D48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E48").Value = 4
$ws.Range("F48").Value = 100112030
$ws.Range("G48").Value = "Poroto granado"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 700
$ws.Range("K48").Value = 23000
$ws.Range("L48").Value = 25000
$ws.Range("M48").Value = 24000
$ws.Range("N48").Value = "`$/malla 25 kilos"
$ws.Range("O48").Value = "Provincia de Limarí"
$ws.Range("P48").Value = 960
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"
